$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (res_partner -> Sheet1)
$ws.Name = "Sheet1"

# ---- Header row (row 1) ----
# Columns A:I already carry the bold header style (s=1); just change the text.
$ws.Range("A1").Value = "Active"
$ws.Range("B1").Value = "City"
$ws.Range("C1").Value = "Country"
$ws.Range("D1").Value = "Date of birth"
$ws.Range("E1").Value = "Display Name"
$ws.Range("F1").Value = "Email"
$ws.Range("G1").Value = "Gender"
$ws.Range("H1").Value = "Is a Company"
$ws.Range("I1").Value = "Job Position"

# New header columns J:M - bold them explicitly to reuse the same bold style.
$ws.Range("J1").Value = "National identification number"
$ws.Range("J1").Font.Bold = $true
$ws.Range("K1").Value = "Phone"
$ws.Range("K1").Font.Bold = $true
$ws.Range("L1").Value = "Related Company"
$ws.Range("L1").Font.Bold = $true
$ws.Range("M1").Value = "Salesperson"
$ws.Range("M1").Font.Bold = $true

# ---- Data row (row 2) ----
# A2: Active = TRUE (was a text name before)
$ws.Range("A2").Value = $true

# B2: City = Bruxelles (was a date value with a date-format style before)
$ws.Range("B2").ClearFormats()
$ws.Range("B2").Value = "Bruxelles"
$ws.Range("B2").WrapText = $true

# C2: Country = empty (was already empty)
$ws.Range("C2").WrapText = $true

# D2: Date of birth = 1998-04-04 (serial 35889), keep yyyy-mm-dd format
$ws.Range("D2").Value = 35889
$ws.Range("D2").NumberFormat = "yyyy\-mm\-dd"

# E2: Display Name = Jhon Doe (was already empty)
$ws.Range("E2").Value = "Jhon Doe"

# F2: Email = empty (previously held "Male" text - must clear it)
$ws.Range("F2").ClearContents()
$ws.Range("F2").WrapText = $true

# G2: Gender = Male (previously held "Watermael-Boitsfort")
$ws.Range("G2").Value = "Male"

# H2: Is a Company = FALSE (previously held "95041328785" text)
$ws.Range("H2").Value = $false

# I2: Job Position = empty (was already empty)
$ws.Range("I2").WrapText = $true

# J2: National identification number = 78586315785 (brand-new cell)
$ws.Range("J2").Value = 78586315785
$ws.Range("J2").WrapText = $true

# K2, L2, M2: Phone / Related Company / Salesperson = empty (brand-new cells)
$ws.Range("K2").WrapText = $true
$ws.Range("L2").WrapText = $true
$ws.Range("M2").WrapText = $true

# Selection matches the saved view (B2 active)
$null = $ws.Range("B2").Select()
